$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All pipelines")

# Insert two new rows before the current row 4 ("Olet" row), shifting it down to row 6.
$ws.Rows.Item(4).Resize(2).Insert()

# The newly inserted rows (4:5) should contain a copy of the original rows 2:3
# (the "FVF"/"FVR" detail rows), before those rows' A-column values are renamed.
$ws.Range("A2:P3").Copy()
$ws.Range("A4").PasteSpecial()
$excel.CutCopyMode = $false

# Rename the values in the original rows 2 and 3 (column A) to the new codes.
$ws.Cells.Item(2, 1).Value = "VF"
$ws.Cells.Item(3, 1).Value = "VR"

$ws.Range("A4").Select()
